$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(" ZAVALETA   ")

# 1) Drop the constant "Nombre-Cliente" column (old column C) - everything right of it shifts left.
$ws.Columns.Item(3).Delete()

# 2) Insert a new title row at the very top.
$ws.Rows.Item(1).Insert()

# 3) Title cell: client name, bold/large, left aligned.
$titleCell = $ws.Range("B1")
$titleCell.Value = "(1140)..ABASTOS 4 CARNES ZAVALETA"
$titleCell.Font.Bold = $true
$titleCell.Font.Size = 18
$titleCell.HorizontalAlignment = -4131
$ws.Rows.Item(1).RowHeight = 23.25

# 4) Bold the data columns A:C (now holding Fecha / Folio / Importe).
$ws.Range("A2:C33").Font.Bold = $true

# 5) Center column B.
$ws.Columns.Item(2).HorizontalAlignment = -4108

# 6) Fix the total formula (it now lives in C, one row further down).
$ws.Range("C35").Formula = "=SUM(C2:C34)"

$ws.Range("L13").Select()
